# Update Typhon_Profits market-board profit figures per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 442.22223
$ws.Range("I33").Value = 442.22223
$ws.Range("K33").Value = 442.22223
$ws.Range("M33").Value = -213.22223
$ws.Range("H70").Value = 1515.4166
$ws.Range("I70").Value = 1397.8572
$ws.Range("J70").Value = 1680
$ws.Range("K70").Value = 4193.571599999999
$ws.Range("L70").Value = 5040
$ws.Range("M70").Value = -3923.571599999999
$ws.Range("N70").Value = -5580
$ws.Range("H73").Value = 1515.4166
$ws.Range("I73").Value = 1397.8572
$ws.Range("J73").Value = 1680
$ws.Range("K73").Value = 4193.571599999999
$ws.Range("L73").Value = 5040
$ws.Range("M73").Value = -3257.571599999999
$ws.Range("N73").Value = -6912
$ws.Range("H98").Value = 1038.75
$ws.Range("I98").Value = 739.8
$ws.Range("J98").Value = 1537
$ws.Range("K98").Value = 739.8
$ws.Range("L98").Value = 1537
$ws.Range("M98").Value = 758.2
$ws.Range("N98").Value = -4533
$ws.Range("H112").Value = 1057
$ws.Range("J112").Value = 1107.9333
$ws.Range("L112").Value = 3323.7999
$ws.Range("N112").Value = -5539.7999
$ws.Range("H116").Value = 4893
$ws.Range("I116").Value = 2416.5
$ws.Range("K116").Value = 2416.5
$ws.Range("M116").Value = 1025.5
$ws.Range("H122").Value = 1038.75
$ws.Range("I122").Value = 739.8
$ws.Range("J122").Value = 1537
$ws.Range("K122").Value = 2219.4
$ws.Range("L122").Value = 4611
$ws.Range("M122").Value = 230.6000000000004
$ws.Range("N122").Value = -9511
$ws.Range("H125").Value = 210.0625
$ws.Range("I125").Value = 249.66667
$ws.Range("J125").Value = 159.14285
$ws.Range("K125").Value = 2247.00003
$ws.Range("L125").Value = 1432.28565
$ws.Range("M125").Value = 212.9999699999998
$ws.Range("N125").Value = -6352.28565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 692.4231
$ws.Range("I2").Value = 788.3889
$ws.Range("J2").Value = 476.5
$ws.Range("K2").Value = 788.3889
$ws.Range("L2").Value = 476.5
$ws.Range("M2").Value = -675.3889
$ws.Range("N2").Value = -702.5
$ws.Range("H32").Value = 7441.015
$ws.Range("I32").Value = 5330.1875
$ws.Range("J32").Value = 13069.889
$ws.Range("K32").Value = 5330.1875
$ws.Range("L32").Value = 13069.889
$ws.Range("M32").Value = -5043.1875
$ws.Range("N32").Value = -13643.889
$ws.Range("H45").Value = 2278.84
$ws.Range("I45").Value = 1600.4375
$ws.Range("J45").Value = 3484.889
$ws.Range("K45").Value = 1600.4375
$ws.Range("L45").Value = 3484.889
$ws.Range("M45").Value = -1223.4375
$ws.Range("N45").Value = -4238.889
$ws.Range("H74").Value = 23810762
$ws.Range("I74").Value = 32258456
$ws.Range("J74").Value = 3619.3635
$ws.Range("K74").Value = 32258456
$ws.Range("L74").Value = 3619.3635
$ws.Range("M74").Value = -32257582
$ws.Range("N74").Value = -5367.363499999999
$ws.Range("H77").Value = 23810762
$ws.Range("I77").Value = 32258456
$ws.Range("J77").Value = 3619.3635
$ws.Range("K77").Value = 161292280
$ws.Range("L77").Value = 18096.8175
$ws.Range("M77").Value = -161287912
$ws.Range("N77").Value = -26832.8175
$ws.Range("H114").Value = 31316.166
$ws.Range("J114").Value = 31316.166
$ws.Range("L114").Value = 31316.166
$ws.Range("N114").Value = -39994.166
$ws.Range("H116").Value = 692.4231
$ws.Range("I116").Value = 788.3889
$ws.Range("J116").Value = 476.5
$ws.Range("K116").Value = 788.3889
$ws.Range("L116").Value = 476.5
$ws.Range("M116").Value = 1505.6111
$ws.Range("N116").Value = -5064.5
$ws.Range("H122").Value = 3148.2222
$ws.Range("I122").Value = 2291.875
$ws.Range("K122").Value = 6875.625
$ws.Range("M122").Value = -4425.625
$ws.Range("H140").Value = 56809.668
$ws.Range("J140").Value = 56809.668
$ws.Range("L140").Value = 56809.668
$ws.Range("N140").Value = -67169.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 692.4231
$ws.Range("I3").Value = 788.3889
$ws.Range("J3").Value = 476.5
$ws.Range("K3").Value = 788.3889
$ws.Range("L3").Value = 476.5
$ws.Range("M3").Value = -674.3889
$ws.Range("N3").Value = -704.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 355
$ws.Range("I7").Value = 103
$ws.Range("J7").Value = 463
$ws.Range("K7").Value = 103
$ws.Range("L7").Value = 463
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = -689
$ws.Range("H31").Value = 4844.913
$ws.Range("I31").Value = 8941.333000000001
$ws.Range("K31").Value = 8941.333000000001
$ws.Range("M31").Value = -8646.333000000001
$ws.Range("H34").Value = 4844.913
$ws.Range("I34").Value = 8941.333000000001
$ws.Range("K34").Value = 8941.333000000001
$ws.Range("M34").Value = -8739.333000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3219
$ws.Range("I3").Value = 1473.75
$ws.Range("J3").Value = 10200
$ws.Range("K3").Value = 4421.25
$ws.Range("L3").Value = 30600
$ws.Range("M3").Value = -4309.25
$ws.Range("N3").Value = -30824
$ws.Range("H113").Value = 816.8570999999999
$ws.Range("I113").Value = 626.3333
$ws.Range("J113").Value = 959.75
$ws.Range("K113").Value = 1878.9999
$ws.Range("L113").Value = 2879.25
$ws.Range("M113").Value = 291.0001
$ws.Range("N113").Value = -7219.25
$ws.Range("H125").Value = 5000
$ws.Range("J125").Value = 5000
$ws.Range("L125").Value = 15000
$ws.Range("N125").Value = -24840
$ws.Range("H129").Value = 284199.38
$ws.Range("I129").Value = 851.6667
$ws.Range("J129").Value = 425873.25
$ws.Range("K129").Value = 2555.0001
$ws.Range("L129").Value = 1277619.75
$ws.Range("M129").Value = 2444.9999
$ws.Range("N129").Value = -1287619.75
$ws.Range("H131").Value = 728.71
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 733.0404
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 2199.1212
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -12279.1212

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3419.087
$ws.Range("I126").Value = 2319.9412
$ws.Range("K126").Value = 6959.823600000001
$ws.Range("M126").Value = -4489.823600000001
$ws.Range("H135").Value = 57585
$ws.Range("J135").Value = 57585
$ws.Range("L135").Value = 57585
$ws.Range("N135").Value = -67725

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 13639
$ws.Range("I48").Value = 11458.5
$ws.Range("K48").Value = 11458.5
$ws.Range("M48").Value = -10797.5
$ws.Range("H135").Value = 28143
$ws.Range("J135").Value = 28143
$ws.Range("L135").Value = 28143
$ws.Range("N135").Value = -38283

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4611.1113
$ws.Range("I62").Value = 3166.6667
$ws.Range("J62").Value = 5333.3335
$ws.Range("K62").Value = 3166.6667
$ws.Range("L62").Value = 5333.3335
$ws.Range("M62").Value = -2542.6667
$ws.Range("N62").Value = -6581.3335
$ws.Range("H65").Value = 4611.1113
$ws.Range("I65").Value = 3166.6667
$ws.Range("J65").Value = 5333.3335
$ws.Range("K65").Value = 15833.3335
$ws.Range("L65").Value = 26666.6675
$ws.Range("M65").Value = -12713.3335
$ws.Range("N65").Value = -32906.6675
$ws.Range("H107").Value = 65341530
$ws.Range("I107").Value = 125000184
$ws.Range("K107").Value = 375000552
$ws.Range("M107").Value = -374998632
$ws.Range("H126").Value = 1362.5927
$ws.Range("I126").Value = 1218.619
$ws.Range("K126").Value = 3655.857
$ws.Range("M126").Value = -1185.857
$ws.Range("J137").Value = 45619.168
$ws.Range("L137").Value = 45619.168
$ws.Range("N137").Value = -55819.168
